# SkillCardData.xlsx edit: rename the "手牌基数" (hand-size-base) upgrade
# card to "手指" (Finger) and update its effect text so the bonus is "+1"
# per copy (instead of "+2"), while still granting "+1 hand-size-base".
#
# Everything else in the sheet (row 2-5, 7-16) keeps the same visible text;
# the many shared-string index churn seen in the target diff is just a
# side effect of the string table being re-packed once the old "手牌基数"
# strings are no longer referenced - that happens automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Card name (column A, row 6)
$ws.Range("A6").Value = "手指"

# Card effect text (column D, row 6) - note the embedded line break matches
# the original formatting ("<br>" followed by a newline before the second
# sentence).
$ws.Range("D6").Value = "每拥有1张《手指》，本牌的点数加1。<br>`n手牌基数加1。"

# The row auto-sizes for wrapped text when the real app re-flows it; pin it
# to the height the content settles at after the edit.
$ws.Rows.Item(6).RowHeight = 99.75

# Update the active selection/view to match where the editor ended up.
$ws.Range("F7").Select()
